$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update column C ("Förändrad") for all existing data rows (2-488):
#    every value moves from 45202 to 45203.
$ws.Range("C2:C488").Value = 45203

# 2. Row 488 picks up an explicit custom row height (matches every other
#    data row in the sheet, which already carries ht="15" customHeight="1").
$ws.Rows.Item(488).RowHeight = 15

# 3. Append the new record as row 489.
$ws.Cells.Item(489, 1).Value = "A 47350-2023"
$ws.Cells.Item(489, 2).Value = 45202
$ws.Cells.Item(489, 3).Value = 45203
$ws.Cells.Item(489, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(489, 5).Value = "EKSJÖ"
$ws.Cells.Item(489, 6).Value = "Sveaskog"
$ws.Cells.Item(489, 7).Value = 1.8
$ws.Cells.Item(489, 8).Value = 0
$ws.Cells.Item(489, 9).Value = 0
$ws.Cells.Item(489, 10).Value = 0
$ws.Cells.Item(489, 11).Value = 0
$ws.Cells.Item(489, 12).Value = 0
$ws.Cells.Item(489, 13).Value = 0
$ws.Cells.Item(489, 14).Value = 0
$ws.Cells.Item(489, 15).Value = 0
$ws.Cells.Item(489, 16).Value = 0
$ws.Cells.Item(489, 17).Value = 0

# Match the date-formatted style used by columns B/C elsewhere in the sheet.
$ws.Cells.Item(489, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(489, 3).NumberFormat = "YYYY-MM-DD"

# Column R is always present (even when empty) with the wrap-text style.
$ws.Cells.Item(489, 18).WrapText = $true

Write-Output "edit applied"
